# Adds the "Przypadek_idealny" (ideal case) worksheet as the first tab,
# reproducing the nominal-metrics data + the abbreviation legend table,
# and tidies up the sheet selections left behind on METIS_AINFO_1m.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no arguments drops the new sheet immediately before
# the currently active sheet - i.e. right at the front of the tab strip,
# which is where "Przypadek_idealny" needs to land.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Przypadek_idealny"

$src = $wb.Worksheets.Item("METIS_AINFO_1m")

# Copy formatting only (no values) from the existing sheet so the new cells
# reuse the very same style entries (bold header / bordered legend table)
# instead of Excel fabricating brand-new style records.
$src.Range("A1:G1").Copy()
$newSheet.Range("A1:G1").PasteSpecial(-4122)
$src.Range("K1:M1").Copy()
$newSheet.Range("J1:L1").PasteSpecial(-4122)
$src.Range("K2:M5").Copy()
$newSheet.Range("J2:L5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row (columns A-G) + legend header (columns J-L)
$newSheet.Range("A1").Value = "Pattern"
$newSheet.Range("B1").Value = "MSE"
$newSheet.Range("C1").Value = "MAE"
$newSheet.Range("D1").Value = "SSD"
$newSheet.Range("E1").Value = "SSD Znormalizowane"
$newSheet.Range("F1").Value = "CC"
$newSheet.Range("G1").Value = "CC Znormalizowane"
$newSheet.Range("J1").Value = "Skrót"
$newSheet.Range("K1").Value = "Rozwinięcie"
$newSheet.Range("L1").Value = "Opis"

# Legend table body (rows 2-5)
$newSheet.Range("J2").Value = "MSE"
$newSheet.Range("K2").Value = "Mean Square Error"
$newSheet.Range("L2").Value = "Im mniejsza wartość MSE, tym obrazy są bardziej podobne."

$newSheet.Range("J3").Value = "MAE"
$newSheet.Range("K3").Value = "Mean Absolute Error"
$newSheet.Range("L3").Value = "Ocenia, jaka jest przeciętna różnica absolutna między odpowiadającymi sobie komórkami danych."

$newSheet.Range("J4").Value = "SSD"
$newSheet.Range("K4").Value = "Sum Square Difference"
$newSheet.Range("L4").Value = "Im mniejsza wartość SSD, tym większe podobieństwo obrazów. Wrażliwy na duże różnice w wartościach (podnoszone do kwadratu)."

$newSheet.Range("J5").Value = "CC"
$newSheet.Range("K5").Value = "Cross - Correlation"
$newSheet.Range("L5").Value = "Miara podobieństwa między obrazami, uwzględniająca przesunięcie (lub dopasowanie) jednego względem drugiego. Im większa wartość tym bardziej dopasowane."

# Data rows 2-28, columns A-G (Cells.Item avoids array-to-range assignment,
# which this host does not support).
$newSheet.Cells.Item(2, 1).Value = 1
$newSheet.Cells.Item(2, 2).Value = 0
$newSheet.Cells.Item(2, 3).Value = 0
$newSheet.Cells.Item(2, 4).Value = 0
$newSheet.Cells.Item(2, 5).Value = 0
$newSheet.Cells.Item(2, 6).Value = 9607126294.2763863
$newSheet.Cells.Item(2, 7).Value = 5276.1340691461583

$newSheet.Cells.Item(3, 1).Value = 2
$newSheet.Cells.Item(3, 2).Value = 0
$newSheet.Cells.Item(3, 3).Value = 0
$newSheet.Cells.Item(3, 4).Value = 0
$newSheet.Cells.Item(3, 5).Value = 0
$newSheet.Cells.Item(3, 6).Value = 9537078371.3914661
$newSheet.Cells.Item(3, 7).Value = 5251.9100268661396

$newSheet.Cells.Item(4, 1).Value = 3
$newSheet.Cells.Item(4, 2).Value = 0
$newSheet.Cells.Item(4, 3).Value = 0
$newSheet.Cells.Item(4, 4).Value = 0
$newSheet.Cells.Item(4, 5).Value = 0
$newSheet.Cells.Item(4, 6).Value = 9704967193.7886658
$newSheet.Cells.Item(4, 7).Value = 5308.7961694239802

$newSheet.Cells.Item(5, 1).Value = 4
$newSheet.Cells.Item(5, 2).Value = 0
$newSheet.Cells.Item(5, 3).Value = 0
$newSheet.Cells.Item(5, 4).Value = 0
$newSheet.Cells.Item(5, 5).Value = 0
$newSheet.Cells.Item(5, 6).Value = 9655891525.2495346
$newSheet.Cells.Item(5, 7).Value = 5291.4859472277149

$newSheet.Cells.Item(6, 1).Value = 5
$newSheet.Cells.Item(6, 2).Value = 0
$newSheet.Cells.Item(6, 3).Value = 0
$newSheet.Cells.Item(6, 4).Value = 0
$newSheet.Cells.Item(6, 5).Value = 0
$newSheet.Cells.Item(6, 6).Value = 9493180688.7898598
$newSheet.Cells.Item(6, 7).Value = 5243.6565160137106

$newSheet.Cells.Item(7, 1).Value = 6
$newSheet.Cells.Item(7, 2).Value = 0
$newSheet.Cells.Item(7, 3).Value = 0
$newSheet.Cells.Item(7, 4).Value = 0
$newSheet.Cells.Item(7, 5).Value = 0
$newSheet.Cells.Item(7, 6).Value = 9508391068.4208641
$newSheet.Cells.Item(7, 7).Value = 5246.0045394441267

$newSheet.Cells.Item(8, 1).Value = 7
$newSheet.Cells.Item(8, 2).Value = 0
$newSheet.Cells.Item(8, 3).Value = 0
$newSheet.Cells.Item(8, 4).Value = 0
$newSheet.Cells.Item(8, 5).Value = 0
$newSheet.Cells.Item(8, 6).Value = 9341253730.1831856
$newSheet.Cells.Item(8, 7).Value = 5191.5504412444143

$newSheet.Cells.Item(9, 1).Value = 8
$newSheet.Cells.Item(9, 2).Value = 0
$newSheet.Cells.Item(9, 3).Value = 0
$newSheet.Cells.Item(9, 4).Value = 0
$newSheet.Cells.Item(9, 5).Value = 0
$newSheet.Cells.Item(9, 6).Value = 9308171782.690506
$newSheet.Cells.Item(9, 7).Value = 5202.2256621001252

$newSheet.Cells.Item(10, 1).Value = 9
$newSheet.Cells.Item(10, 2).Value = 0
$newSheet.Cells.Item(10, 3).Value = 0
$newSheet.Cells.Item(10, 4).Value = 0
$newSheet.Cells.Item(10, 5).Value = 0
$newSheet.Cells.Item(10, 6).Value = 9708761142.989624
$newSheet.Cells.Item(10, 7).Value = 5304.2327583285196

$newSheet.Cells.Item(11, 1).Value = 10
$newSheet.Cells.Item(11, 2).Value = 0
$newSheet.Cells.Item(11, 3).Value = 0
$newSheet.Cells.Item(11, 4).Value = 0
$newSheet.Cells.Item(11, 5).Value = 0
$newSheet.Cells.Item(11, 6).Value = 9631875498.8010292
$newSheet.Cells.Item(11, 7).Value = 5268.5634248686883

$newSheet.Cells.Item(12, 1).Value = 11
$newSheet.Cells.Item(12, 2).Value = 0
$newSheet.Cells.Item(12, 3).Value = 0
$newSheet.Cells.Item(12, 4).Value = 0
$newSheet.Cells.Item(12, 5).Value = 0
$newSheet.Cells.Item(12, 6).Value = 9435064987.509409
$newSheet.Cells.Item(12, 7).Value = 5210.9746656120251

$newSheet.Cells.Item(13, 1).Value = 12
$newSheet.Cells.Item(13, 2).Value = 0
$newSheet.Cells.Item(13, 3).Value = 0
$newSheet.Cells.Item(13, 4).Value = 0
$newSheet.Cells.Item(13, 5).Value = 0
$newSheet.Cells.Item(13, 6).Value = 9459449891.3005981
$newSheet.Cells.Item(13, 7).Value = 5220.3136188947292

$newSheet.Cells.Item(14, 1).Value = 13
$newSheet.Cells.Item(14, 2).Value = 0
$newSheet.Cells.Item(14, 3).Value = 0
$newSheet.Cells.Item(14, 4).Value = 0
$newSheet.Cells.Item(14, 5).Value = 0
$newSheet.Cells.Item(14, 6).Value = 9474660075.6223602
$newSheet.Cells.Item(14, 7).Value = 5226.2520645193972

$newSheet.Cells.Item(15, 1).Value = 14
$newSheet.Cells.Item(15, 2).Value = 0
$newSheet.Cells.Item(15, 3).Value = 0
$newSheet.Cells.Item(15, 4).Value = 0
$newSheet.Cells.Item(15, 5).Value = 0
$newSheet.Cells.Item(15, 6).Value = 9412365552.6067162
$newSheet.Cells.Item(15, 7).Value = 5205.4834851355899

$newSheet.Cells.Item(16, 1).Value = 15
$newSheet.Cells.Item(16, 2).Value = 0
$newSheet.Cells.Item(16, 3).Value = 0
$newSheet.Cells.Item(16, 4).Value = 0
$newSheet.Cells.Item(16, 5).Value = 0
$newSheet.Cells.Item(16, 6).Value = 9390075622.0031052
$newSheet.Cells.Item(16, 7).Value = 5232.4967643419232

$newSheet.Cells.Item(17, 1).Value = 16
$newSheet.Cells.Item(17, 2).Value = 0
$newSheet.Cells.Item(17, 3).Value = 0
$newSheet.Cells.Item(17, 4).Value = 0
$newSheet.Cells.Item(17, 5).Value = 0
$newSheet.Cells.Item(17, 6).Value = 9311995833.2055168
$newSheet.Cells.Item(17, 7).Value = 5185.6943767462699

$newSheet.Cells.Item(18, 1).Value = 17
$newSheet.Cells.Item(18, 2).Value = 0
$newSheet.Cells.Item(18, 3).Value = 0
$newSheet.Cells.Item(18, 4).Value = 0
$newSheet.Cells.Item(18, 5).Value = 0
$newSheet.Cells.Item(18, 6).Value = 9246922578.4275494
$newSheet.Cells.Item(18, 7).Value = 5174.9807500982988

$newSheet.Cells.Item(19, 1).Value = 18
$newSheet.Cells.Item(19, 2).Value = 0
$newSheet.Cells.Item(19, 3).Value = 0
$newSheet.Cells.Item(19, 4).Value = 0
$newSheet.Cells.Item(19, 5).Value = 0
$newSheet.Cells.Item(19, 6).Value = 9622749855.0244923
$newSheet.Cells.Item(19, 7).Value = 5291.7073133219919

$newSheet.Cells.Item(20, 1).Value = 19
$newSheet.Cells.Item(20, 2).Value = 0
$newSheet.Cells.Item(20, 3).Value = 0
$newSheet.Cells.Item(20, 4).Value = 0
$newSheet.Cells.Item(20, 5).Value = 0
$newSheet.Cells.Item(20, 6).Value = 9537654980.4943905
$newSheet.Cells.Item(20, 7).Value = 5266.0777400080378

$newSheet.Cells.Item(21, 1).Value = 20
$newSheet.Cells.Item(21, 2).Value = 0
$newSheet.Cells.Item(21, 3).Value = 0
$newSheet.Cells.Item(21, 4).Value = 0
$newSheet.Cells.Item(21, 5).Value = 0
$newSheet.Cells.Item(21, 6).Value = 9350827020.8372917
$newSheet.Cells.Item(21, 7).Value = 5204.4192387542071

$newSheet.Cells.Item(22, 1).Value = 21
$newSheet.Cells.Item(22, 2).Value = 0
$newSheet.Cells.Item(22, 3).Value = 0
$newSheet.Cells.Item(22, 4).Value = 0
$newSheet.Cells.Item(22, 5).Value = 0
$newSheet.Cells.Item(22, 6).Value = 9543680685.8393612
$newSheet.Cells.Item(22, 7).Value = 5247.9698237055063

$newSheet.Cells.Item(23, 1).Value = 22
$newSheet.Cells.Item(23, 2).Value = 0
$newSheet.Cells.Item(23, 3).Value = 0
$newSheet.Cells.Item(23, 4).Value = 0
$newSheet.Cells.Item(23, 5).Value = 0
$newSheet.Cells.Item(23, 6).Value = 9327255860.2523842
$newSheet.Cells.Item(23, 7).Value = 5193.0167709651842

$newSheet.Cells.Item(24, 1).Value = 23
$newSheet.Cells.Item(24, 2).Value = 0
$newSheet.Cells.Item(24, 3).Value = 0
$newSheet.Cells.Item(24, 4).Value = 0
$newSheet.Cells.Item(24, 5).Value = 0
$newSheet.Cells.Item(24, 6).Value = 9226484504.3392639
$newSheet.Cells.Item(24, 7).Value = 5160.5891535931523

$newSheet.Cells.Item(25, 1).Value = 24
$newSheet.Cells.Item(25, 2).Value = 0
$newSheet.Cells.Item(25, 3).Value = 0
$newSheet.Cells.Item(25, 4).Value = 0
$newSheet.Cells.Item(25, 5).Value = 0
$newSheet.Cells.Item(25, 6).Value = 9507893480.7103882
$newSheet.Cells.Item(25, 7).Value = 5256.9630666262756

$newSheet.Cells.Item(26, 1).Value = 25
$newSheet.Cells.Item(26, 2).Value = 0
$newSheet.Cells.Item(26, 3).Value = 0
$newSheet.Cells.Item(26, 4).Value = 0
$newSheet.Cells.Item(26, 5).Value = 0
$newSheet.Cells.Item(26, 6).Value = 9359400377.6798515
$newSheet.Cells.Item(26, 7).Value = 5203.8795893225979

$newSheet.Cells.Item(27, 1).Value = 26
$newSheet.Cells.Item(27, 2).Value = 0
$newSheet.Cells.Item(27, 3).Value = 0
$newSheet.Cells.Item(27, 4).Value = 0
$newSheet.Cells.Item(27, 5).Value = 0
$newSheet.Cells.Item(27, 6).Value = 9182893492.4142342
$newSheet.Cells.Item(27, 7).Value = 5146.3367599380499

$newSheet.Cells.Item(28, 1).Value = 27
$newSheet.Cells.Item(28, 2).Value = 0
$newSheet.Cells.Item(28, 3).Value = 0
$newSheet.Cells.Item(28, 4).Value = 0
$newSheet.Cells.Item(28, 5).Value = 0
$newSheet.Cells.Item(28, 6).Value = 9115839822.3403645
$newSheet.Cells.Item(28, 7).Value = 5126.7030433348427

# Column widths to roughly match the bestFit widths used elsewhere in the
# workbook (column E / G / L). The host quantizes ColumnWidth to whole
# pixels, so these land close to, not bit-exact with, the source widths.
$newSheet.Columns.Item(5).ColumnWidth = 16.33
$newSheet.Columns.Item(7).ColumnWidth = 15.5
$newSheet.Columns.Item(12).ColumnWidth = 125.5

# METIS_AINFO_1m is no longer the active tab, but it keeps a new selection
# (K1:M5 instead of the old K21).
$src.Range("K1:M5").Select()

# "Przypadek_idealny" ends up the active/selected tab with I9 selected.
$newSheet.Select()
$newSheet.Range("I9").Select()
